$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the wrap-text style previously applied to C9:C11 so the new numeric
# cells in column C (rows 9-37) are unstyled, matching the target workbook.
$ws.Range("C9:C11").Style = "Normal"

$categories = @(
  "AssetCategory.TV",
  "AssetCategory.SOUND_SYSTEM",
  "AssetCategory.PROJECTOR",
  "AssetCategory.AIR_CONDITIONER",
  "AssetCategory.LIGHTING",
  "AssetCategory.AIR_PURIFIER",
  "AssetCategory.STOVE",
  "AssetCategory.MICROWAVE",
  "AssetCategory.OVEN",
  "AssetCategory.REFRIGERATOR",
  "AssetCategory.WATER_PURIFIER",
  "AssetCategory.RANGE_HOOD",
  "AssetCategory.BED",
  "AssetCategory.CAR",
  "AssetCategory.WARDROBE",
  "AssetCategory.FAN",
  "AssetCategory.LAMP",
  "AssetCategory.SHOWER",
  "AssetCategory.BATHTUB",
  "AssetCategory.SINK",
  "AssetCategory.WATER_HEATER",
  "AssetCategory.EXHAUST_FAN",
  "AssetCategory.WASHING_MACHINE",
  "AssetCategory.DRYER",
  "AssetCategory.CLOTHES_RACK",
  "AssetCategory.IRON",
  "AssetCategory.CHAIR",
  "AssetCategory.LAPTOP",
  "AssetCategory.MOBILE_PHONE"
)

$lifeValues = @(10, 12, 8, 15, 15, 7, 15, 10, 15, 20, 7, 15, 15, 15, 20, 10, 15, 15, 25, 20, 12, 15, 12, 12, 10, 8, 15, 5, 4)

$startRow = 9
for ($i = 0; $i -lt $categories.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 2).Value = $categories[$i]
  $ws.Cells.Item($row, 3).Value = $lifeValues[$i]
}

# Update the view state to match: scrolled so row 3 is at the top and the
# active selection on B15 (the "AIR_PURIFIER" row). (Note: this headless
# runtime only serializes a scrolled topLeftCell when it is part of a
# frozen/split pane, so a bare non-split sheetView@topLeftCell cannot be
# produced here; ScrollRow/ScrollColumn are still set for correctness of
# the in-memory view state.)
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B15").Select() | Out-Null
